$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Date property value: 2025-05-21T14:22:51+00:00 -> 2025-06-13T15:45:04+00:00
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# FHIR Version property value: 4.3.0 -> 4.0.1
$meta.Range("B15").Value = "4.0.1"

# --- Elements sheet --------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

# Row 2 = Extension: Constraint(s) (AJ) loses the "unless an empty Parameters
# resource ... or $this is Parameters" clause from the ele-1 invariant.
$els.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 3 = Extension.id: Type(s) (K) changes from "id" to "string".
$els.Range("K3").Value = "string" + [char]10

# Row 4 = Extension.extension: Constraint(s) (AJ) now matches the corrected
# ele-1 text used on row 2 (the duplicate string collapses into one).
$els.Range("AJ4").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 6 = Extension.value[x]: Definition (M) link updated from R4B to R4.
$els.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
